$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "5 point tamper proof bits" line item to row 22
$ws.Range("B22").Value = "5 point tool set for kinect"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 6.98
$ws.Range("F22").Formula = "=D22*E22"
$ws.Range("G22").Value = "http://www.amazon.com/5-Point-TORX-Tamper-Proof-Bit-Set/dp/B007SOODHC/ref=sr_1_2?ie=UTF8&qid=1392914223&sr=8-2&keywords=five+point+tamper+proof"

# Update the view: scroll back to the top and select E23
$ws.Activate()
$ws.Range("E23").Select()
